$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as text so that numeric-looking strings
    # (e.g. "1.00", "0.999") are preserved exactly instead of being
    # coerced into a number, then restore the default ("Normal") style so
    # no stray number-format style is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "64.878.11"
$ws.Range("E2").Value = "  +3.86%  "
$ws.Range("D3").Value = "3.102.10"
$ws.Range("E3").Value = "  +2.24%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("D5") "561.56"
$ws.Range("E5").Value = "  +3.18%  "
Set-TextValue $ws.Range("D6") "144.28"
$ws.Range("E6").Value = "  +7.31%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.099.30"
$ws.Range("E8").Value = "  +2.44%  "
Set-TextValue $ws.Range("D9") "0.501"
$ws.Range("E9").Value = "  +1.25%  "
Set-TextValue $ws.Range("D10") "6.33"
$ws.Range("E10").Value = "  +3.00%  "
Set-TextValue $ws.Range("D11") "0.153"
$ws.Range("E11").Value = "  +3.42%  "
Set-TextValue $ws.Range("D12") "0.474"
$ws.Range("E12").Value = "  +6.05%  "
Set-TextValue $ws.Range("D13") "0.0000230"
$ws.Range("E13").Value = "  +3.13%  "
Set-TextValue $ws.Range("D14") "35.47"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "3.600.06"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").Value = "64.859.53"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("D17").Value = "3.091.85"
$ws.Range("E17").Value = "  +1.92%  "
Set-TextValue $ws.Range("D18") "0.110"
$ws.Range("E18").Value = "  +1.35%  "
Set-TextValue $ws.Range("D19") "6.80"
$ws.Range("E19").Value = "  +2.48%  "
Set-TextValue $ws.Range("D20") "481.35"
$ws.Range("E20").Value = "  +0.51%  "
Set-TextValue $ws.Range("D21") "13.81"
$ws.Range("E21").Value = "  +3.93%  "
Set-TextValue $ws.Range("D22") "0.690"
$ws.Range("E22").Value = "  +2.23%  "
Set-TextValue $ws.Range("D23") "7.60"
$ws.Range("E23").Value = "  +7.81%  "
Set-TextValue $ws.Range("D24") "13.55"
$ws.Range("E24").Value = "  +11.23%  "
Set-TextValue $ws.Range("D25") "81.40"
$ws.Range("E25").Value = "  +0.50%  "
Set-TextValue $ws.Range("D26") "1.00"
$ws.Range("E26").Value = "  +0.05%  "
Set-TextValue $ws.Range("D27") "2.79"
$ws.Range("E27").Value = "  +2.73%  "
Set-TextValue $ws.Range("D28") "8.26"
$ws.Range("E28").Value = "  +5.58%  "
Set-TextValue $ws.Range("D29") "2.07"
$ws.Range("E29").Value = "  +6.73%  "
Set-TextValue $ws.Range("D30") "0.996"
$ws.Range("E30").Value = "  -0.08%  "
Set-TextValue $ws.Range("D31") "26.21"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("E32").Value = "  +2.17%  "
Set-TextValue $ws.Range("D33") "2.49"
$ws.Range("E33").Value = "  +4.69%  "
Set-TextValue $ws.Range("D34") "5.64"
$ws.Range("E34").Value = "  -0.53%  "
Set-TextValue $ws.Range("D35") "6.18"
$ws.Range("E35").Value = "  +5.04%  "
Set-TextValue $ws.Range("D36") "54.91"
$ws.Range("E36").Value = "  +0.19%  "
Set-TextValue $ws.Range("D37") "470.01"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D38") "3.01"
$ws.Range("E38").Value = "  +21.11%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.0839"
$ws.Range("E39").Value = "  +4.46%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.0412"
$ws.Range("E40").Value = "  +5.63%  "
$ws.Range("D41").Value = "2.980.24"
$ws.Range("E41").Value = "  -5.91%  "
Set-TextValue $ws.Range("D42") "8.28"
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("E43").Value = "  -2.63%  "
Set-TextValue $ws.Range("D44") "28.32"
$ws.Range("E44").Value = "  +6.71%  "
Set-TextValue $ws.Range("D45") "0.261"
$ws.Range("E45").Value = "  +6.36%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D46") "2.17"
$ws.Range("E46").Value = "  +9.21%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D47") "1.00"
$ws.Range("E47").Value = "  +0.00%  "
Set-TextValue $ws.Range("D48") "0.113"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("D49").Value = "0.0₃0528"
$ws.Range("E49").Value = "  +4.91%  "
Set-TextValue $ws.Range("D50") "116.85"
$ws.Range("E50").Value = "  +1.69%  "
Set-TextValue $ws.Range("D51") "2.08"
$ws.Range("E51").Value = "  +2.71%  "
